$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D and E columns (rows 2-51) to prevent Excel from
# auto-converting numeric-looking strings into numbers/percentages.
$fmtRange = $ws.Range("D2:E51")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.518.88"
$ws.Range("E2").Value = "  +3.16%  "
$ws.Range("D3").Value = "3.374.61"
$ws.Range("E3").Value = "  +4.72%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "192.56"
$ws.Range("E5").Value = "  +5.52%  "
$ws.Range("D6").Value = "594.08"
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("E9").Value = "  +3.44%  "
$ws.Range("D10").Value = "6.76"
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("D11").Value = "0.422"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").Value = "3.963.17"
$ws.Range("E12").Value = "  +4.85%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("E14").Value = "  +3.91%  "
$ws.Range("D15").Value = "69.555.93"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "3.389.47"
$ws.Range("E17").Value = "  +6.29%  "
$ws.Range("D18").Value = "450.72"
$ws.Range("E18").Value = "  +14.08%  "
$ws.Range("D19").Value = "5.84"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("D20").Value = "13.79"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("D21").Value = "7.83"
$ws.Range("E21").Value = "  +3.80%  "
$ws.Range("D22").Value = "73.48"
$ws.Range("E22").Value = "  +3.89%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "3.519.01"
$ws.Range("E24").Value = "  +4.60%  "
$ws.Range("D25").Value = "0.519"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +4.00%  "
$ws.Range("E27").Value = "  +4.62%  "
$ws.Range("D28").Value = "9.60"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "2.00"
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").Value = "23.25"
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("E33").Value = "  +5.14%  "
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +3.96%  "
$ws.Range("D37").Value = "164.69"
$ws.Range("E37").Value = "  +2.40%  "
$ws.Range("D38").Value = "1.94"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("D39").Value = "27.24"
$ws.Range("E39").Value = "  +3.85%  "
$ws.Range("D40").Value = "0.823"
$ws.Range("E40").Value = "  +2.55%  "
$ws.Range("E41").Value = "  +1.21%  "
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").Value = "2.746.00"
$ws.Range("E43").Value = "  +5.62%  "
$ws.Range("D44").Value = "2.53"
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("D45").Value = "25.51"
$ws.Range("E45").Value = "  +4.28%  "
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").Value = "344.83"
$ws.Range("E47").Value = "  +3.61%  "
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("E49").Value = "  +3.64%  "
$ws.Range("D50").Value = "32.92"
$ws.Range("E50").Value = "  +7.33%  "
$ws.Range("D51").Value = "1.03"
$ws.Range("E51").Value = "  +7.77%  "

# Restore the original default cell style (removes the temporary text
# number format override) while keeping the values as text.
$fmtRange.Style = "Normal"
